$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by exactly one day.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}

# Updated "Notified Production (MW)" readings (column B) for the rows whose
# values changed as part of the NRG location adjustment.
$bUpdates = @{
    23 = 2
    24 = 3
    25 = 3
    26 = 8
    27 = 9
    28 = 11
    29 = 21
    30 = 136
    31 = 190
    32 = 254
    33 = 311
    34 = 584
    35 = 676
    36 = 785
    37 = 884
    38 = 1138
    39 = 1232
    40 = 1333
    41 = 1421
    42 = 1586
    43 = 1653
    44 = 1713
    45 = 1762
    46 = 1794
    47 = 1812
    48 = 1821
    49 = 1816
    50 = 1792
    51 = 1771
    52 = 1743
    53 = 1711
    54 = 1616
    55 = 1569
    56 = 1513
    57 = 1460
    58 = 1293
    59 = 1221
    60 = 1142
    61 = 1059
    62 = 848
    63 = 758
    64 = 674
    65 = 586
    66 = 377
    67 = 295
    68 = 226
    69 = 172
    70 = 52
    71 = 22
    72 = 15
    73 = 14
    74 = 2
    75 = 2
    76 = 1
    77 = 1
    78 = 1
    79 = 1
    80 = 1
    81 = 1
}

foreach ($row in $bUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value2 = $bUpdates[$row]
}
